$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H ("Save") header - copy the same header formatting (bold,
# centered, bordered) used by the other header cells in row 1, then set
# the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" stat (0/1) for every data row, 2 through 23.
$values = @(0, 0, 0, 0, 1, 1, 0, 0, 0, 1, 1, 1, 1, 0, 0, 1, 0, 1, 0, 1, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
